# Applies the edits described by the commit diff:
#  1. "Objetivos" paragraph: change how the app notifies users (text -> sound/lights).
#  2. "Funcionais" requirement: drop "notificação do app(som)" in favor of "luz".
#  3. "Não-Funcionais" requirement: drop the duplicated "luz" and "tema no app".
#  4. "Exigências Legais" requirement: drop the trailing "app com fácil entendimento".

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "por meio de um app os responsáveis da limpeza a recolher esse lixo. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "por meio de onda sonora em um painel, por luzinhas quando está cheia.",
    2) | Out-Null

$d.Content.Find.Execute(
    " programação, notificação do app(som).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " programação, luz",
    2) | Out-Null

$d.Content.Find.Execute(
    "Decoração, luz, tamanho do sensor, tema no app.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Decoração, tamanho do sensor.",
    2) | Out-Null

$d.Content.Find.Execute(
    " a passo para não ter problemas com o resultado do projeto, app com fácil entendimento.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " a passo para não ter problemas com o resultado do projeto.",
    2) | Out-Null
